$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Facebook/Twitter count values (column B) for rows 2-12
$ws.Range("B2").Value = 1904
$ws.Range("B3").Value = 1929
$ws.Range("B4").Value = 1935
$ws.Range("B5").Value = 1927
$ws.Range("B6").Value = 1941
$ws.Range("B7").Value = 1877
$ws.Range("B8").Value = 1871
$ws.Range("B9").Value = 1876
$ws.Range("B10").Value = 1834
$ws.Range("B11").Value = 1933
$ws.Range("B12").Value = 1936

# Add three new quarterly rows (13-15), copying formatting from the last
# existing data row so the date style (numFmt) carries over correctly.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C15").PasteSpecial(-4122)

$ws.Range("A13").Value = 44561
$ws.Range("B13").Value = 1934
$ws.Range("C13").Value = "Q"

$ws.Range("A14").Value = 44651
$ws.Range("B14").Value = 1956
$ws.Range("C14").Value = "Q"

$ws.Range("A15").Value = 44742
$ws.Range("B15").Value = 1993
$ws.Range("C15").Value = "Q"
